$sheetsToAdd = @(
    @{
        Name = '04Aug2020 23.24.06'
        Rows = @(
            ,@('https://supchainplanmerch-mcp-user-interface-photonqa.services.west.nonprod.wsgc.com/')
            ,@('Page Name', 'DOM Load Finish Time', 'Page Interaction Time')
            ,@('Login Page', 3059, 54246)
        )
    }
    @{
        Name = '05Aug2020 05.43.05'
        Rows = @(
            ,@('https://supchainplanmerch-mcp-user-interface-photonqa.services.west.nonprod.wsgc.com/')
            ,@('Page Name', 'DOM Load Finish Time', 'Page Interaction Time')
            ,@('Login Page', 1096, 12893)
            ,@('Libraries > Country', 745, 3580)
            ,@('Login Page', 1253, 10131)
            ,@('Libraries > Vendor', 686, 6077)
            ,@('Login Page', 1079, 8503)
            ,@('Libraries > Image', 863, 6603)
            ,@('Login Page', 840, 8482)
            ,@('Libraries > Collection', 4780, 20628)
            ,@('Login Page', 1304, 30796)
            ,@('Libraries > Lading Port', 1229, 11694)
            ,@('Login Page', 5223, 37333)
            ,@('Libraries > Agent', 2737, 11592)
            ,@('Login Page', 2200, 44180)
            ,@('Libraries > Season Collection', 1281, 5976)
            ,@('Login Page', 2554, 51017)
            ,@('Libraries > Diffs', 1828, 16645)
            ,@('Login Page', 2517, 43379)
            ,@('Libraries > JDA Hierarchy', 2245, 20001)
            ,@('Login Page', 7594, 42823)
            ,@('Libraries > Business Objects', 1705, 15053)
            ,@('Login Page', 6245, 53084)
            ,@('Libraries > Item', 1378, 23072)
            ,@('Login Page', 4388, 43313)
            ,@('Libraries > Item-Pack-Components', 3440, 5293)
            ,@('Login Page', 3698, 33810)
            ,@('Libraries > Item-Season', 2533, 16231)
            ,@('Login Page', 4282, 32499)
            ,@('Libraries > Item-Parent-Season', 6000, 16261)
            ,@('Login Page', 1665, 14196)
            ,@('Libraries > Merch Hierarchy', 12111, 3535)
            ,@('Login Page', 1073, 16605)
            ,@('Libraries > Season', 668, 10)
            ,@('Login Page', 751, 7404)
            ,@('Login Page', 2028, 37823)
            ,@('Libraries > Item-Source-Season', 1743, 15846)
            ,@('Login Page', 2112, 34629)
            ,@('Libraries > Merch List', 2405, 28515)
        )
    }
    @{
        Name = '05Aug2020 06.12.05'
        Rows = @(
            ,@('https://supchainplanmerch-mcp-user-interface-photonqa.services.west.nonprod.wsgc.com/')
            ,@('Page Name', 'DOM Load Finish Time', 'Page Interaction Time')
            ,@('Login Page', 2309, 43951)
            ,@('Libraries > Item-Parent', 4341, 14223)
        )
    }
    @{
        Name = '09Aug2020 16.30.28'
        Rows = @(
            ,@('https://supchainplanmerch-mcp-user-interface-photonqa.services.west.nonprod.wsgc.com/')
            ,@('Page Name', 'DOM Load Finish Time', 'Page Interaction Time')
            ,@('Login Page', 1160, 13399)
            ,@('User Profile', 1160, 642)
            ,@('Libraries > Merch List- Page 1', 741, 12324)
            ,@('Libraries > Merch List- Page 2', 741, 8817)
            ,@('Libraries > Merch List- Page 3', 741, 15221)
            ,@('Libraries > Merch List- Page 4', 741, 14338)
            ,@('Libraries > Merch List- Page 5', 741, 22694)
            ,@('User Profile', 741, 720)
            ,@('Libraries > Merch List- Page 1', 636, 13084)
            ,@('Libraries > Merch List- Page 2', 636, 10239)
            ,@('Libraries > Merch List- Page 3', 636, 14596)
            ,@('Libraries > Merch List- Page 4', 636, 14924)
            ,@('Libraries > Merch List- Page 5', 636, 21389)
        )
    }
)

$wb = $excel.ActiveWorkbook

foreach ($sheetDef in $sheetsToAdd) {
    $lastIndex = $wb.Worksheets.Count
    $lastSheet = $wb.Worksheets.Item($lastIndex)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $sheetDef.Name

    $rowNum = 1
    foreach ($row in $sheetDef.Rows) {
        $colNum = 1
        foreach ($val in $row) {
            if ($null -ne $val) {
                $ws.Cells.Item($rowNum, $colNum).Value = $val
            }
            $colNum = $colNum + 1
        }
        $rowNum = $rowNum + 1
    }
}
